# excel_merge_neighbors_3.py used group function which works perfectly
#
# The sheet holds a small task table in A1:G7. Row 4 ("pb2"/"p2-t1") had its
# numeric values sitting two columns to the right (F4:G4) instead of lining
# up with its sibling rows (D:E) - a artifact of grouping rows that don't
# all share the same populated columns. The fix re-aligns row 4 by moving
# F4:G4 into D4:E4, then (as the neighbor-merge / group step does for every
# row) pads the table out with the same header/empty-cell pattern all the
# way to column P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-align row 4: its values were offset by two columns (F,G -> D,E) ---
$ws.Range("D4").Value = $ws.Range("F4").Value2
$ws.Range("E4").Value = $ws.Range("G4").Value2
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""

# --- Extend the header row (row 1) from G1 out to P1, continuing 0..15 ---
# (Writing these also grows the sheet's used range/dimension from G7 to P7;
#  the rest of the table's newly-in-range cells stay implicitly blank.)
$headerCols = @("H","I","J","K","L","M","N","O","P")
$headerVals = @(7,8,9,10,11,12,13,14,15)
$ws.Range("A1").Copy()
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $cell = $ws.Range($headerCols[$i] + "1")
    $cell.Value = $headerVals[$i]
    $cell.PasteSpecial(-4122)   # xlPasteFormats - match A1's bordered/bold style
}
$excel.CutCopyMode = $false
